$wb = $excel.ActiveWorkbook

$layout = $wb.Worksheets.Item("Complete_Run_Layout")
$layout.Activate()

$layout.Range("C4").Value = 2
$layout.Range("C5").Value = 2

$layout.Range("H2").Select()
